$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = 'Other found locations'

# Row 2 (Zhou et al. / PMC7196181): update authors string and add source column
$ws.Range("E2").Value = '[Shuang-Jiang%Zhou%NULL%0,    Li-Gang%Zhang%NULL%0,    Lei-Lei%Wang%NULL%0,    Zhao-Chang%Guo%NULL%0,    Jing-Qi%Wang%NULL%0,    Jin-Cheng%Chen%NULL%0,    Mei%Liu%NULL%0,    Xi%Chen%NULL%0,    Jing-Xu%Chen%chenjx1110@163.com%0]'
$ws.Range("I2").Value = '_PMC_Springer'

# Row 3 (Oosterhoff et al. / PMC7205689): update authors string and add source column
$ws.Range("E3").Value = '[Benjamin%Oosterhoff%Benjamin.oosterhoff@montana.edu%0,    Cara A.%Palmer%NULL%1,    Jenna%Wilson%NULL%1,    Natalie%Shook%NULL%1]'
$ws.Range("I3").Value = '_PMC_elsevier'

# Row 4 (Secer et al. / PMC7293436): update authors string and add source column
$ws.Range("E4").Value = '[İsmail%Seçer%ismailsecer84@gmail.com%0,    Sümeyye%Ulaş%NULL%2,    Sümeyye%Ulaş%NULL%0]'
$ws.Range("I4").Value = '_PMC_Springer'

# Row 5 (Qu et al. / PMC9281280): update authors string and add source column
$ws.Range("E5").Value = '[Miao%Qu%NULL%1,    Kun%Yang%NULL%0,    Yujia%Cao%NULL%1,    Mei Hong%Xiu%xiumeihong97@163.com%1,    Xiang Yang%Zhang%zhangxy@psych.ac.cn%2,    Xiang Yang%Zhang%zhangxy@psych.ac.cn%0]'
$ws.Range("I5").Value = '_PMC_Springer'
